$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.562.00'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '3.394.98'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''559.39'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').Value = '''175.81'
$ws.Range('D7').Value = '''0.629'
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('D8').Value = '3.386.31'
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = '''0.174'
$ws.Range('E10').Value = '  +5.16%  '
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').Value = '''53.86'
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('D14').Value = '''9.23'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').Value = '3.930.97'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '''18.39'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '''0.119'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.380.08'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').Value = '65.362.97'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').Value = '''457.81'
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').Value = '''4.90'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Value = '''4.11'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').Value = '''14.09'
$ws.Range('E25').Value = '  +5.24%  '
$ws.Range('D26').Value = '''87.50'
$ws.Range('E26').Value = '  +1.46%  '
$ws.Range('D27').Value = '''2.92'
$ws.Range('E27').Value = '  +2.68%  '
$ws.Range('D28').Value = '''10.72'
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('D29').Value = '''8.77'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Value = '''31.23'
$ws.Range('E30').Value = '  +3.82%  '
$ws.Range('D31').Value = '''6.57'
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').Value = '''63.58'
$ws.Range('E32').Value = '  +7.04%  '
$ws.Range('D33').Value = '''11.52'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('D34').Value = '''583.88'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').Value = '''3.60'
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('D39').Value = '''35.88'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').Value = '''0.375'
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').Value = '0.0₃0746'
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').Value = '3.097.98'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('D44').Value = '''2.79'
$ws.Range('E44').Value = '  -1.85%  '
$ws.Range('D45').Value = '''2.45'
$ws.Range('E45').Value = '  -2.66%  '
$ws.Range('E46').Value = '  +2.15%  '
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('D48').Value = '''0.999'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '''140.59'
$ws.Range('E49').Value = '  +3.63%  '
$ws.Range('D50').Value = '''2.53'
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').Value = '''3.12'
$ws.Range('E51').Value = '  +8.73%  '
